# Update the HALO sub-assembly/harness BOM list on Sheet1.
# - Columns: A=Title, E=Description, F=Project, G=Material, H=SurfaceFinish,
#            I=Commodity, O=UnitOfMeasure, P=MakeOrBuy,
#            Q=DrawnBy, S=CheckedBy, U=EngApproval, W=MfgApproval, Y=QAApproval
# The previous placeholder/test rows (a/d/qe/f/.../qeqeqe) are replaced with
# real part data, row 4 & 6 get a "Material <not specified>" remark (with
# trailing CR/LF artifacts exactly as produced by SolidWorks' export), and a
# new row 7 (LS5 / HALO) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $title, $desc, $material) {
    $ws.Cells.Item($r, 1).Value  = $title       # A - Title
    $ws.Cells.Item($r, 5).Value  = $desc         # E - Description
    $ws.Cells.Item($r, 6).Value  = "HALO"        # F - Project
    $ws.Cells.Item($r, 7).Value  = $material     # G - Material
    $ws.Cells.Item($r, 8).Value  = "--"          # H - SurfaceFinish
    $ws.Cells.Item($r, 9).Value  = "CBL"         # I - Commodity
    $ws.Cells.Item($r, 15).Value = "PC"          # O - UnitOfMeasure
    $ws.Cells.Item($r, 16).Value = "MAKE"        # P - MakeOrBuy
    $ws.Cells.Item($r, 17).Value = "--"          # Q - DrawnBy
    $ws.Cells.Item($r, 19).Value = "--"          # S - CheckedBy
    $ws.Cells.Item($r, 21).Value = "--"          # U - EngApproval
    $ws.Cells.Item($r, 23).Value = "--"          # W - MfgApproval
    $ws.Cells.Item($r, 25).Value = "--"          # Y - QAApproval
}

$matPlain  = "Material <not specified>"
$matCrLf   = "Material <not specified>`r`n"
$matCrLf2  = "Material <not specified>`r`n`r`n"

Set-Row 2 "LS1"  "Bumper Shell Interface" $matPlain
Set-Row 3 "LS2"  "Base-Shell Breakout (From Main Body Components to Shell Baseplate)" $matPlain
Set-Row 4 "LS3"  "Rear Aux (Kernel Module to Rear Fan)" $matCrLf
Set-Row 5 "LS4a" "NVR (Kernel Module to NVR)" $matPlain
Set-Row 6 "LS4b" "Front Aux (From Kernel Module to Speakers + Intake Fan)`r`n" $matCrLf2
Set-Row 7 "LS5"  "Kernel-Hat Link (From Top Enclosure Interface Plate to Kernel Module)" $matCrLf2

# Columns J:N (Manufacturer, ManufacturerPartNo, Supplier, SupplierPartNo,
# WebLink) are no longer populated for these rows - remove the leftover
# placeholder values from the old test data.
$ws.Range("J2:N7").ClearContents() | Out-Null

# Let the columns whose content actually changed width re-autofit (as Excel
# does automatically after the data entry), same columns the workbook
# session ended up resizing.
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(10).AutoFit() | Out-Null

$ws.Application.ActiveWindow.Zoom = 129
$ws.Range("E13").Select() | Out-Null
